$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers: build an OOXML run sequence and inject it into a paragraph via
# Range.InsertXML, which (unlike Find.Execute / Range.Text) preserves each
# run as a distinct <w:r> element instead of merging adjacent identically
# formatted runs together.
# ---------------------------------------------------------------------------

function Build-RunsXml($runsInfo) {
    $sb = New-Object System.Text.StringBuilder
    foreach ($info in $runsInfo) {
        $text = $info[0]
        $italic = $info[1]
        $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
        if ($italic) {
            [void]$sb.Append('<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>')
        }
        else {
            [void]$sb.Append('<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>')
        }
    }
    return $sb.ToString()
}

function Set-ParagraphRuns($paragraph, $runsInfo) {
    $r = $paragraph.Range
    # Exclude the trailing paragraph-mark character from the replaced range.
    $r.SetRange($r.Start, $r.End - 1)
    $r.Text = ""
    $ins = $d.Range($r.Start, $r.Start)
    $runsXml = Build-RunsXml $runsInfo
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
    $ins.InsertXML($pkg)
}

function Find-ParagraphIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

$LDQ = [string][char]0x201C   # “
$RDQ = [string][char]0x201D   # ”

# ---------------------------------------------------------------------------
# 1) "Frame the prediction of an unknown quantity using an *interval*." ->
#    "Summarize the likely range of a variable using an *interval*."
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex("Frame the prediction of an unknown quantity")
$runs = @(
    , @("Summarize the likely range of a variable using an", $false)
    , @(" ", $false)
    , @("interval", $true)
    , @(".", $false)
)
Set-ParagraphRuns $d.Paragraphs.Item($idx) $runs

# ---------------------------------------------------------------------------
# 2) "Choose the prediction interval endpoints to convey a range of likely
#    outcomes." ->
#    "Choose the summary interval endpoints to include the "vast majority"
#    of the values."
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex("Choose the prediction interval endpoints")
$runs = @(
    , @("Choose the summary interval endpoints to include the", $false)
    , @(" ", $false)
    , @($LDQ, $false)
    , @("vast majority", $false)
    , @($RDQ, $false)
    , @(" ", $false)
    , @("of the values.", $false)
)
Set-ParagraphRuns $d.Paragraphs.Item($idx) $runs

# ---------------------------------------------------------------------------
# 3) "Understand that predictions are more likely to be right when the
#    interval is constructed to be long." ->
#    "Master the convention that "vast majority" is the central 95% of the
#    values."
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex("Understand that predictions are more likely to be right")
$runs = @(
    , @("Master the convention that", $false)
    , @(" ", $false)
    , @($LDQ, $false)
    , @("vast majority", $false)
    , @($RDQ, $false)
    , @(" ", $false)
    , @("is the central 95% of the values.", $false)
)
Set-ParagraphRuns $d.Paragraphs.Item($idx) $runs

# ---------------------------------------------------------------------------
# 4) Delete the paragraph: "Understand, in contrast to (3), that predictions
#    are more likely to be informative when the interval is short."
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex("Understand, in contrast to \(3\)")
$d.Paragraphs.Item($idx).Range.Delete()

# ---------------------------------------------------------------------------
# 5) Delete the paragraph: "Master the convention that is widely used for
#    accomplishing both (3) and (4), e.g. the use of 95% coverage as a
#    sensible evaluation of the trade-off. Identify 95% as a "level"
#    selected specifically to have a reasonable balance between (3) and (4)."
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex("widely used for accomplishing both")
$d.Paragraphs.Item($idx).Range.Delete()

# ---------------------------------------------------------------------------
# 6) "Be able to work with two formats for describing an interval: "A to B"
#    and "C ± D"" ->
#    "Be able to work with two formats for describing an interval: "A to B"
#    and "C ± D"."  (append a trailing period)
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex("Be able to work with two formats")
$runs = @(
    , @("Be able to work with two formats for describing an interval:", $false)
    , @(" ", $false)
    , @($LDQ, $false)
    , @("A to B", $false)
    , @($RDQ, $false)
    , @(" ", $false)
    , @("and", $false)
    , @(" ", $false)
    , @($LDQ, $false)
    , @("C ± D", $false)
    , @($RDQ, $false)
    , @(".", $false)
)
Set-ParagraphRuns $d.Paragraphs.Item($idx) $runs

Write-Host "Edits applied successfully."
